$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire column A, shifting B:F left to A:E.
$ws.Range("A:A").Delete()
